# Adjust Investment Summary table column widths for better formatting.
#
# Three tables (slides 2, 3, 4 -- "Why This Solution?", "Business Value -
# Financial Impact" and "Risk Mitigation") each get their right-most grid
# column nudged by 1 EMU (1/12700 pt), which also nudges the overall table
# extent by 1 EMU. The placeholder/sample copy that filled each cell is
# cleared out as part of the same formatting pass.

$p = $ppt.ActivePresentation

function Clear-TableCells($tbl) {
    $rows = $tbl.Rows.Count
    $cols = $tbl.Columns.Count
    for ($r = 1; $r -le $rows; $r++) {
        for ($c = 1; $c -le $cols; $c++) {
            $tbl.Cell($r, $c).Shape.TextFrame.TextRange.Text = ""
        }
    }
}

# --- Slide 2: "Why This Solution?" (2-column table) -----------------------
$s2 = $p.Slides.Item(2)
$tbl2 = $s2.Shapes.Item("Table Placeholder 3").Table
$tbl2.Columns.Item(2).Width = 4355467 / 12700.0
Clear-TableCells $tbl2

# --- Slide 3: "Business Value - Financial Impact" (2-column table) --------
$s3 = $p.Slides.Item(3)
$tbl3 = $s3.Shapes.Item("Table Placeholder 3").Table
$tbl3.Columns.Item(2).Width = 4355467 / 12700.0
Clear-TableCells $tbl3

# --- Slide 4: "Risk Mitigation" (3-column table) ---------------------------
$s4 = $p.Slides.Item(4)
$tbl4 = $s4.Shapes.Item("Table Placeholder 3").Table
$tbl4.Columns.Item(3).Width = 2903645 / 12700.0
Clear-TableCells $tbl4
